$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general_settings")
$ws.Range("B7").Value = "Prod02"
